$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "bad drivers" sample counts / roaming % for the weekly refresh
$ws.Range("C3").Value = 55
$ws.Range("D3").Value = 98.40000000000001
$ws.Range("C4").Value = 55

# Narrow the spacer/value columns to match this week's layout
$ws.Columns.Item(2).ColumnWidth = 13.17
$ws.Columns.Item(5).ColumnWidth = 1.17

# No good drivers this week - replace the header row with a single message
# and drop the now-unused blank rows below it.
$ws.Range("A11:E11").Clear()
$ws.Range("A11").Value = "No good drivers found."
$ws.Rows("12:16").Delete()
